# ------------------------------------------------------------------
# Applies the template's daily regen:
#  1) bump the cached "datetimeFigureOut" field text on the slide
#     master and every slide layout from 11/18/2025 -> 11/19/2025.
#  2) on slide 3, move/resize the existing background rectangle and
#     add a second rectangle (duplicate of the first, same style)
#     positioned/sized per the new layout.
# ------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Date placeholder text bump (slide master + all layouts) ---

$newDate = "11/19/2025"

$masterShapes = $p.SlideMaster.Shapes
for ($j = 1; $j -le $masterShapes.Count; $j++) {
    $sh = $masterShapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $lyt = $layouts.Item($i)
    $lytShapes = $lyt.Shapes
    for ($j = 1; $j -le $lytShapes.Count; $j++) {
        $sh = $lytShapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2) Slide 3 rectangle reflow ---

$s3 = $p.Slides.Item(3)
$rect1 = $s3.Shapes.Item(1)

# Duplicate the existing rectangle before moving it, so the new
# shape inherits the same fill / line / style / text formatting.
$rect2 = $rect1.Duplicate()
$rect2.Name = "Rectangle 10"

# New shape ("Rectangle 10"): off (3737702,250584) ext (4833711,3010702) EMU
$rect2.Left = 3737702 / 12700
$rect2.Top = 250584 / 12700
$rect2.Width = 4833711 / 12700
$rect2.Height = 3010702 / 12700

# Existing shape ("Rectangle 139"): off (2757513,3596714) ext (6794089,3010703) EMU
$rect1.Left = 2757513 / 12700
$rect1.Top = 3596714 / 12700
$rect1.Width = 6794089 / 12700
$rect1.Height = 3010703 / 12700
